$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$vals = @(5.6, 1.9, 6.1, 5.9, 2, 6.3, 5.8, 2.1, 4.4, 1.4)
for ($i=0; $i -lt $vals.Length; $i++) {
  $ws.Cells.Item(160+$i, 1).Value = $vals[$i]
}
